# Auto-generated Excel COM-interop script to apply the diff changes
# Workbook: Asura Profits -- per-row H..N recalculated price/profit values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (ALC, item id 5512)
$ws.Range("H33").Value = 163.75
$ws.Range("I33").Value = 163.77777
$ws.Range("J33").Value = 163.66667
$ws.Range("K33").Value = 163.77777
$ws.Range("L33").Value = 163.66667
$ws.Range("M33").Value = 65.22223
$ws.Range("N33").Value = -621.6666700000001

# Row 64 (ALC, item id 5506)
$ws.Range("H64").Value = 3027.342
$ws.Range("I64").Value = 2708.1765
$ws.Range("J64").Value = 3285.7144
$ws.Range("K64").Value = 2708.1765
$ws.Range("L64").Value = 3285.7144
$ws.Range("M64").Value = -2460.1765
$ws.Range("N64").Value = -3781.7144

# Row 67 (ALC, item id 5506)
$ws.Range("H67").Value = 3027.342
$ws.Range("I67").Value = 2708.1765
$ws.Range("J67").Value = 3285.7144
$ws.Range("K67").Value = 2708.1765
$ws.Range("L67").Value = 3285.7144
$ws.Range("M67").Value = -1850.1765
$ws.Range("N67").Value = -5001.7144

# Row 70 (ALC, item id 12604)
$ws.Range("H70").Value = 63754.375
$ws.Range("J70").Value = 1315.8334
$ws.Range("L70").Value = 3947.5002
$ws.Range("N70").Value = -4487.5002

# Row 73 (ALC, item id 12604)
$ws.Range("H73").Value = 63754.375
$ws.Range("J73").Value = 1315.8334
$ws.Range("L73").Value = 3947.5002
$ws.Range("N73").Value = -5819.5002

# Row 81 (ALC, item id 10637)
$ws.Range("H81").Value = 24100
$ws.Range("J81").Value = 24100
$ws.Range("L81").Value = 24100
$ws.Range("N81").Value = -26096

# Row 84 (ALC, item id 10637)
$ws.Range("H84").Value = 24100
$ws.Range("J84").Value = 24100
$ws.Range("L84").Value = 72300
$ws.Range("N84").Value = -82284

# Row 100 (ALC, item id 19906)
$ws.Range("H100").Value = 3258.611
$ws.Range("J100").Value = 3317
$ws.Range("L100").Value = 3317
$ws.Range("N100").Value = -4399

# Row 137 (ALC, item id 44013)
$ws.Range("H137").Value = 1456.25
$ws.Range("I137").Value = 1215.3334
$ws.Range("J137").Value = 2179
$ws.Range("K137").Value = 3646.0002
$ws.Range("L137").Value = 6537
$ws.Range("M137").Value = -1096.0002
$ws.Range("N137").Value = -11637

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (ARM, item id 27714)
$ws.Range("H45").Value = 1089.7778
$ws.Range("I45").Value = 982.4
$ws.Range("K45").Value = 982.4
$ws.Range("M45").Value = -605.4

# Row 61 (ARM, item id 43999)
$ws.Range("H61").Value = 2597.5833
$ws.Range("I61").Value = 2019
$ws.Range("J61").Value = 4333.3335
$ws.Range("K61").Value = 2019
$ws.Range("L61").Value = 4333.3335
$ws.Range("M61").Value = -1807
$ws.Range("N61").Value = -4757.3335

# Row 123 (ARM, item id 34107)
$ws.Range("H123").Value = 30228
$ws.Range("J123").Value = 30228
$ws.Range("L123").Value = 30228
$ws.Range("N123").Value = -40028

# Row 132 (ARM, item id 43997)
$ws.Range("H132").Value = 3232
$ws.Range("I132").Value = 2334.2
$ws.Range("J132").Value = 4129.8
$ws.Range("K132").Value = 7002.599999999999
$ws.Range("L132").Value = 12389.4
$ws.Range("M132").Value = -4472.599999999999
$ws.Range("N132").Value = -17449.4

# Row 136 (ARM, item id 43999)
$ws.Range("H136").Value = 2597.5833
$ws.Range("I136").Value = 2019
$ws.Range("J136").Value = 4333.3335
$ws.Range("K136").Value = 6057
$ws.Range("L136").Value = 13000.0005
$ws.Range("M136").Value = -3507
$ws.Range("N136").Value = -18100.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (CRP, item id 44021)
$ws.Range("H58").Value = 1816.8572
$ws.Range("I58").Value = 1816.8572
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1816.8572
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1613.8572
$ws.Range("N58").ClearContents()

# Row 136 (CRP, item id 44021)
$ws.Range("H136").Value = 1816.8572
$ws.Range("I136").Value = 1816.8572
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5450.571599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2900.571599999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 75 (CUL, item id 12863)
$ws.Range("H75").Value = 3678.6667
$ws.Range("I75").Value = 904
$ws.Range("J75").Value = 6453.3335
$ws.Range("K75").Value = 2712
$ws.Range("L75").Value = 19360.0005
$ws.Range("M75").Value = -1714
$ws.Range("N75").Value = -21356.0005

# Row 78 (CUL, item id 12863)
$ws.Range("H78").Value = 3678.6667
$ws.Range("I78").Value = 904
$ws.Range("J78").Value = 6453.3335
$ws.Range("K78").Value = 8136
$ws.Range("L78").Value = 58080.0015
$ws.Range("M78").Value = -3144
$ws.Range("N78").Value = -68064.0015

$ws = $wb.Worksheets.Item("GSM")
# Row 20 (GSM, item id 4095)
$ws.Range("H20").Value = 35000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20490

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW, item id 36249)
$ws.Range("H7").Value = 6166.6665
$ws.Range("I7").Value = 6250
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 6250
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -6138
$ws.Range("N7").Value = -6224

# Row 68 (LTW, item id 12563)
$ws.Range("H68").Value = 2067.6924
$ws.Range("I68").Value = 1742.8572
$ws.Range("J68").Value = 2446.6667
$ws.Range("K68").Value = 1742.8572
$ws.Range("L68").Value = 2446.6667
$ws.Range("M68").Value = -993.8571999999999
$ws.Range("N68").Value = -3944.6667

# Row 71 (LTW, item id 12563)
$ws.Range("H71").Value = 2067.6924
$ws.Range("I71").Value = 1742.8572
$ws.Range("J71").Value = 2446.6667
$ws.Range("K71").Value = 8714.286
$ws.Range("L71").Value = 12233.3335
$ws.Range("M71").Value = -4970.286
$ws.Range("N71").Value = -19721.3335

# Row 126 (LTW, item id 36249)
$ws.Range("H126").Value = 6166.6665
$ws.Range("I126").Value = 6250
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 18750
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -16280
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (WVR, item id 3307)
$ws.Range("H2").Value = 1002
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 62 (WVR, item id 12589)
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

# Row 65 (WVR, item id 12589)
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

# Row 69 (WVR, item id 10951)
$ws.Range("H69").Value = 21423.666
$ws.Range("J69").Value = 21423.666
$ws.Range("L69").Value = 21423.666
$ws.Range("N69").Value = -22921.666

# Row 72 (WVR, item id 10951)
$ws.Range("H72").Value = 21423.666
$ws.Range("J72").Value = 21423.666
$ws.Range("L72").Value = 64270.99800000001
$ws.Range("N72").Value = -71758.99800000001

# Row 75 (WVR, item id 11957)
$ws.Range("H75").Value = 200000
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# Row 78 (WVR, item id 11957)
$ws.Range("H78").Value = 200000
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# Row 96 (WVR, item id 19977)
$ws.Range("H96").Value = 1467.6666
$ws.Range("I96").Value = 1467.6666
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1467.6666
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -94.66660000000002
$ws.Range("N96").ClearContents()

# Row 132 (WVR, item id 44029)
$ws.Range("H132").Value = 2550.389
$ws.Range("I132").Value = 1911.5555
$ws.Range("J132").Value = 3189.2222
$ws.Range("K132").Value = 5734.666499999999
$ws.Range("L132").Value = 9567.6666
$ws.Range("M132").Value = -3204.666499999999
$ws.Range("N132").Value = -14627.6666

# Row 136 (WVR, item id 44031)
$ws.Range("H136").Value = 1929.2667
$ws.Range("I136").Value = 1929.2667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5787.800099999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3237.800099999999
$ws.Range("N136").ClearContents()

Write-Output "Applied all Asura_Profits updates"
